# Apply the 'Product Backlogs.docx' edits described in the commit diff.
# All of the changes below are pure run-consolidations: the visible text
# before and after each Find/Replace is identical - only the underlying
# <w:r>/<w:t> run-splitting of that text changes (Word's Find & Replace
# naturally re-emits a matched range as a single run), except for the
# two spots where new words are actually inserted.

$d = $word.ActiveDocument

# 1) "Given that there is a [map][, when the customer views it they can ][find the restaurant.]"
#    -> merge the 4 runs into a single run (text unchanged).
$d.Content.Find.Execute(
    "Given that there is a map, when the customer views it they can find the restaurant.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Given that there is a map, when the customer views it they can find the restaurant.",
    2) | Out-Null

# 2) "As a customer I want to view an about page which includes ... personal statement."
#    -> append " so that I can learn about the restaurant." after "personal statement"
#       (replacing the trailing ".") and merge runs.
$d.Content.Find.Execute(
    "As a customer I want to view an about page which includes what veganism is, a few pictures of example vegan cuisine and the head chef's personal statement.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As a customer I want to view an about page which includes what veganism is, a few pictures of example vegan cuisine and the head chef's personal statement so that I can learn about the restaurant.",
    2) | Out-Null

# 3) "Given that[ staff member][s][ ][are]" -> merge into one run (text unchanged).
$d.Content.Find.Execute(
    "Given that staff members are able to delete reservations",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Given that staff members are able to delete reservations",
    2) | Out-Null

# 4) "[ the table][ can][ become free for other customers.]" -> merge into one run.
$d.Content.Find.Execute(
    "when they receive a call to cancel, the table can become free for other customers.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "when they receive a call to cancel, the table can become free for other customers.",
    2) | Out-Null

# 5) "Add staff ability to [update][ reservations from database]" -> merge into one run.
$d.Content.Find.Execute(
    "Add staff ability to update reservations from database",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Add staff ability to update reservations from database",
    2) | Out-Null
